$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Preserve string storage: if Excel would auto-convert the literal
    # text into a number, force Text format first so it stays a string,
    # matching the source workbook's original (string-typed) cell content.
    if ($value -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $range.NumberFormat = "@"
    }
    $range.Value = $value
}

Set-TextValue $ws.Range("D2") '28.174.05'
$ws.Range("E2").Value = '  -1.69%  '

Set-TextValue $ws.Range("D3") '1.803.07'
$ws.Range("E3").Value = '  +0.31%  '

$ws.Range("E4").Value = '  +0.03%  '

Set-TextValue $ws.Range("D5") '316.42'
$ws.Range("E5").Value = '  +0.94%  '

$ws.Range("E6").Value = '  +0.01%  '

Set-TextValue $ws.Range("D7") '0.5376'
$ws.Range("E7").Value = '  +1.28%  '

Set-TextValue $ws.Range("D8") '0.3779'
$ws.Range("E8").Value = '  +0.13%  '

Set-TextValue $ws.Range("D9") '0.07463'
$ws.Range("E9").Value = '  -1.05%  '

Set-TextValue $ws.Range("D10") '42.03'
$ws.Range("E10").Value = '  -1.06%  '

Set-TextValue $ws.Range("D11") '1.096'
$ws.Range("E11").Value = '  -2.30%  '

Set-TextValue $ws.Range("D12") '1.0000'

Set-TextValue $ws.Range("D13") '6.203'
$ws.Range("E13").Value = '  +0.03%  '

Set-TextValue $ws.Range("D14") '20.51'
$ws.Range("E14").Value = '  -3.22%  '

Set-TextValue $ws.Range("D15") '7.371'
$ws.Range("E15").Value = '  -1.16%  '

Set-TextValue $ws.Range("D16") '1.803.80'
$ws.Range("E16").Value = '  +0.68%  '

Set-TextValue $ws.Range("D17") '89.76'
$ws.Range("E17").Value = '  -0.85%  '

Set-TextValue $ws.Range("D18") '0.00001063'
$ws.Range("E18").Value = '  -0.35%  '

Set-TextValue $ws.Range("D19") '0.06499'
$ws.Range("E19").Value = '  +0.86%  '

$ws.Range("E20").Value = '  +0.05%  '

Set-TextValue $ws.Range("D21") '17.39'
$ws.Range("E21").Value = '  +0.56%  '

Set-TextValue $ws.Range("D22") '5.926'
$ws.Range("E22").Value = '  -0.11%  '

Set-TextValue $ws.Range("D23") '28.200.56'
$ws.Range("E23").Value = '  -1.60%  '

$ws.Range("E24").Value = '  -0.11%  '

Set-TextValue $ws.Range("D25") '2.090'
$ws.Range("E25").Value = '  -0.25%  '

Set-TextValue $ws.Range("D26") '156.10'
$ws.Range("E26").Value = '  -2.92%  '

Set-TextValue $ws.Range("D27") '20.51'
$ws.Range("E27").Value = '  -0.43%  '

Set-TextValue $ws.Range("D28") '2.009.72'
$ws.Range("E28").Value = '  +0.57%  '

Set-TextValue $ws.Range("D29") '2.326'
$ws.Range("E29").Value = '  -3.51%  '

Set-TextValue $ws.Range("D30") '121.87'
$ws.Range("E30").Value = '  -1.44%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D31") '1.125'
$ws.Range("E31").Value = '  -0.60%  '

$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D32") '0.1112'
$ws.Range("E32").Value = '  +8.80%  '

Set-TextValue $ws.Range("D33") '3.659'
$ws.Range("E33").Value = '  -0.12%  '

Set-TextValue $ws.Range("D34") '5.589'
$ws.Range("E34").Value = '  -2.65%  '

Set-TextValue $ws.Range("D35") '0.07056'
$ws.Range("E35").Value = '  +7.23%  '

Set-TextValue $ws.Range("D36") '0.2221'
$ws.Range("E36").Value = '  -3.57%  '

Set-TextValue $ws.Range("D37") '0.02301'
$ws.Range("E37").Value = '  -1.15%  '

Set-TextValue $ws.Range("D38") '5.080'
$ws.Range("E38").Value = '  -0.20%  '

Set-TextValue $ws.Range("D39") '8.450'
$ws.Range("E39").Value = '  -3.83%  '

Set-TextValue $ws.Range("D40") '0.6171'
$ws.Range("E40").Value = '  -2.47%  '

Set-TextValue $ws.Range("D41") '11.12'
$ws.Range("E41").Value = '  -3.56%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D42") '1.177'
$ws.Range("E42").Value = '  -2.50%  '

$ws.Range("B43").Value = 'WEMIXTOKEN'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D43") '1.429'
$ws.Range("E43").Value = '  +2.60%  '

Set-TextValue $ws.Range("D44") '13.45'
$ws.Range("E44").Value = '  -1.27%  '

$ws.Range("E45").Value = '  +0.36%  '

Set-TextValue $ws.Range("D46") '0.5759'
$ws.Range("E46").Value = '  -3.00%  '

Set-TextValue $ws.Range("D47") '125.22'
$ws.Range("E47").Value = '  -0.46%  '

$ws.Range("E48").Value = '  +1.64%  '

Set-TextValue $ws.Range("D49") '1.926'
$ws.Range("E49").Value = '  -2.88%  '

$ws.Range("E50").Value = '  -1.72%  '

Set-TextValue $ws.Range("D51") '71.86'
$ws.Range("E51").Value = '  -1.47%  '
